$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.06"

$ws.Range("E3").Value = "'7.39%"

$ws.Range("D4").Value = "'5.191"
$ws.Range("E4").Value = "'1.42%"

$ws.Range("D5").Value = "'0.05736"
$ws.Range("E5").Value = "'0.88%"

$ws.Range("D6").Value = "'6.546"
$ws.Range("E6").Value = "'0.45%"

$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'3.089"
$ws.Range("E7").Value = "'2.66%"

$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.8587"
$ws.Range("E8").Value = "'4.79%"

$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").Value = "'0.8706"
$ws.Range("E9").Value = "'1.52%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1367"
$ws.Range("E10").Value = "'2.48%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07071"
$ws.Range("E11").Value = "'1.90%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02998"
$ws.Range("E12").Value = "'4.91%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09384"
$ws.Range("E13").Value = "'-0.09%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001539"
$ws.Range("E14").Value = "'0.49%"

$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005998"
$ws.Range("E15").Value = "'0.14%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006026"
$ws.Range("E16").Value = "'-3.04%"

$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D17").Value = "'0.007489"
$ws.Range("E17").Value = "'5,225.15%"

$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.491"
$ws.Range("E18").Value = "'-0.59%"

$ws.Range("D19").Value = "'2.182"
$ws.Range("E19").Value = "'-5.85%"

$ws.Range("E20").Value = "'1.11%"

$ws.Range("D21").Value = "'0.03343"
$ws.Range("E21").Value = "'4.16%"

$ws.Range("D22").Value = "'0.1293"
$ws.Range("E22").Value = "'1.48%"

$ws.Range("D23").Value = "'3.487"
$ws.Range("E23").Value = "'-2.18%"

$ws.Range("D24").Value = "'0.04147"
$ws.Range("E24").Value = "'3.46%"

$ws.Range("E25").Value = "'0.40%"

$ws.Range("D26").Value = "'0.005002"
$ws.Range("E26").Value = "'12.06%"

$ws.Range("D27").Value = "'0.001224"
$ws.Range("E27").Value = "'0.69%"

$ws.Range("D28").Value = "'0.0001210"
$ws.Range("E28").Value = "'2.55%"

$ws.Range("D40").Value = "'0.03755"
$ws.Range("E40").Value = "'1.01%"

$ws.Range("D41").Value = "'0.005758"
$ws.Range("E41").Value = "'-3.67%"

$ws.Range("D42").Value = "'0.1074"
$ws.Range("E42").Value = "'1.48%"

$ws.Range("D43").Value = "'0.002428"
$ws.Range("E43").Value = "'5.59%"

$ws.Range("D44").Value = "'0.009442"
$ws.Range("E44").Value = "'-2.71%"

$ws.Range("D45").Value = "'0.00005255"
$ws.Range("E45").Value = "'2.83%"

$ws.Range("E46").Value = "'0.00%"

$ws.Range("D47").Value = "'0.05698"
$ws.Range("E47").Value = "'-43.58%"

$ws.Range("E48").Value = "'-9.30%"

$ws.Range("E49").Value = "'0.00%"

$ws.Range("E50").Value = "'0.00%"
